$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 260. This shifts the
# existing rows 260-300 down to 262-302 (Excel relative-reference semantics),
# which is exactly what the target workbook needs since every row from the
# old 260 onward reappears two rows lower, and the former last row (300)
# ends up at 302.
$ws.Rows("260:261").Insert()

# New row 260: Cultivar XV region / Primera, Region de Arica y Parinacota
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 44984
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = 100112043
$ws.Cells.Item(260, 7).Value = "Pepino dulce"
$ws.Cells.Item(260, 8).Value = "Cultivar XV región"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 260
$ws.Cells.Item(260, 11).Value = 18000
$ws.Cells.Item(260, 12).Value = 19000
$ws.Cells.Item(260, 13).Value = 18462
$ws.Cells.Item(260, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(260, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(260, 16).Value = 1026
$ws.Cells.Item(260, 17).Value = 18
$ws.Cells.Item(260, 18).Value = "Hortaliza"

# New row 261: Cultivar XV region / Segunda, Region de Arica y Parinacota
$ws.Cells.Item(261, 1).Value = 10
$ws.Cells.Item(261, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(261, 3).Value = "La Araucanía"
$ws.Cells.Item(261, 4).Value = 44984
$ws.Cells.Item(261, 5).Value = 9
$ws.Cells.Item(261, 6).Value = 100112043
$ws.Cells.Item(261, 7).Value = "Pepino dulce"
$ws.Cells.Item(261, 8).Value = "Cultivar XV región"
$ws.Cells.Item(261, 9).Value = "Segunda"
$ws.Cells.Item(261, 10).Value = 50
$ws.Cells.Item(261, 11).Value = 16000
$ws.Cells.Item(261, 12).Value = 16000
$ws.Cells.Item(261, 13).Value = 16000
$ws.Cells.Item(261, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(261, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(261, 16).Value = 889
$ws.Cells.Item(261, 17).Value = 18
$ws.Cells.Item(261, 18).Value = "Hortaliza"

# Make sure the date cells keep the same date-time number format the rest of
# column D uses (style index carries over from the insert, but set explicitly
# to be safe).
$ws.Range("D260:D261").NumberFormat = $ws.Range("D262").NumberFormat
